# mise a jour des methodes delete de client identity claimobject et unité
#
# The "objects.xlsx" reference sheet documents the fields of the "object de
# plainte" model. This edit clarifies the `severity_level` field description
# (row 1, column E) by appending the list of accepted values, and leaves the
# active selection parked on that same cell (E1) as it was when the sheet
# was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the allowed values to the "severity_level" column header comment.
$ws.Range("E1").Value = "le status du niveau de gravité         Valeurs possibles :(low,medium,high)"

# Leave the selection on the cell that was just edited (E1 only, not the
# former A1:G2 header block).
$ws.Range("E1").Select()
